# Pipeliner Integration - Level 2 - init
#
# The "Hardware Country Code" (R3) and "Software Country Code" (S3) cells
# on the Account sheet both contained the stray value "CH". Clear them
# back to the sheet's blank placeholder value (a single space, matching
# every other untouched cell in row 3) so the now-unused "CH" shared
# string is dropped on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("R3").Value = " "
$ws.Range("S3").Value = " "
